# Refresh market-derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across all eight crafting-sheet tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with the latest pulled market data. Generated by the scheduled runner from the
# authoritative price snapshot; values below are the verbatim new readings.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 414230.88
$ws.Range("J17").Value = 441779.6
$ws.Range("L17").Value = 1325338.8
$ws.Range("N17").Value = -1325674.8
$ws.Range("H41").Value = 1235.125
$ws.Range("I41").Value = 740.0909
$ws.Range("J41").Value = 1654
$ws.Range("K41").Value = 740.0909
$ws.Range("L41").Value = 1654
$ws.Range("M41").Value = -300.0909
$ws.Range("N41").Value = -2534
$ws.Range("H112").Value = 4491.579
$ws.Range("J112").Value = 5036
$ws.Range("L112").Value = 15108
$ws.Range("N112").Value = -17324
$ws.Range("H119").Value = 750
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 750
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 2250
$ws.Range("N119").Value = -11926
$ws.Range("H138").Value = 4973.025
$ws.Range("J138").Value = 4214.961
$ws.Range("L138").Value = 12644.883
$ws.Range("N138").Value = -22924.883
$ws.Range("M119").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1522.9584
$ws.Range("I2").Value = 1023.3125
$ws.Range("J2").Value = 2522.25
$ws.Range("K2").Value = 1023.3125
$ws.Range("L2").Value = 2522.25
$ws.Range("M2").Value = -910.3125
$ws.Range("N2").Value = -2748.25
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("H32").Value = 229709.05
$ws.Range("I32").Value = 278849.94
$ws.Range("J32").Value = 8575
$ws.Range("K32").Value = 278849.94
$ws.Range("L32").Value = 8575
$ws.Range("M32").Value = -278562.94
$ws.Range("N32").Value = -9149
$ws.Range("H61").Value = 2678570.8
$ws.Range("I61").Value = 5273.9443
$ws.Range("J61").Value = 14708407
$ws.Range("K61").Value = 5273.9443
$ws.Range("L61").Value = 14708407
$ws.Range("M61").Value = -5061.9443
$ws.Range("N61").Value = -14708831
$ws.Range("H110").Value = 1297.3334
$ws.Range("I110").Value = 971
$ws.Range("K110").Value = 971
$ws.Range("M110").Value = 1074
$ws.Range("H116").Value = 1522.9584
$ws.Range("I116").Value = 1023.3125
$ws.Range("J116").Value = 2522.25
$ws.Range("K116").Value = 1023.3125
$ws.Range("L116").Value = 2522.25
$ws.Range("M116").Value = 1270.6875
$ws.Range("N116").Value = -7110.25
$ws.Range("H124").Value = 12404.667
$ws.Range("J124").Value = 12404.667
$ws.Range("L124").Value = 12404.667
$ws.Range("N124").Value = -22224.667
$ws.Range("H125").Value = 150000
$ws.Range("J125").Value = 150000
$ws.Range("L125").Value = 150000
$ws.Range("N125").Value = -159840
$ws.Range("H132").Value = 5993.5
$ws.Range("J132").Value = 5995
$ws.Range("L132").Value = 17985
$ws.Range("N132").Value = -23045
$ws.Range("H136").Value = 2678570.8
$ws.Range("I136").Value = 5273.9443
$ws.Range("J136").Value = 14708407
$ws.Range("K136").Value = 15821.8329
$ws.Range("L136").Value = 44125221
$ws.Range("M136").Value = -13271.8329
$ws.Range("N136").Value = -44130321
$ws.Range("H141").Value = 77748.25
$ws.Range("J141").Value = 77748.25
$ws.Range("L141").Value = 77748.25
$ws.Range("N141").Value = -88108.25
$ws.Range("N7").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1522.9584
$ws.Range("I3").Value = 1023.3125
$ws.Range("J3").Value = 2522.25
$ws.Range("K3").Value = 1023.3125
$ws.Range("L3").Value = 2522.25
$ws.Range("M3").Value = -909.3125
$ws.Range("N3").Value = -2750.25
$ws.Range("H20").Value = 953.7857
$ws.Range("J20").Value = 1099.1666
$ws.Range("L20").Value = 1099.1666
$ws.Range("N20").Value = -1593.1666
$ws.Range("H105").Value = 18612.25
$ws.Range("I105").Value = 34966.332
$ws.Range("K105").Value = 34966.332
$ws.Range("M105").Value = -33219.332

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2927.5715
$ws.Range("I2").Value = 3498.25
$ws.Range("K2").Value = 3498.25
$ws.Range("M2").Value = -3385.25
$ws.Range("H33").Value = 6781.75
$ws.Range("I33").Value = 1195.5
$ws.Range("J33").Value = 12368
$ws.Range("K33").Value = 1195.5
$ws.Range("L33").Value = 12368
$ws.Range("M33").Value = -816.5
$ws.Range("N33").Value = -13126
$ws.Range("H58").Value = 2324.9697
$ws.Range("I58").Value = 2000.5
$ws.Range("J58").Value = 3530.1428
$ws.Range("K58").Value = 2000.5
$ws.Range("L58").Value = 3530.1428
$ws.Range("M58").Value = -1797.5
$ws.Range("N58").Value = -3936.1428
$ws.Range("H99").Value = 74287464
$ws.Range("J99").Value = 100001870
$ws.Range("L99").Value = 100001870
$ws.Range("N99").Value = -100004866
$ws.Range("H126").Value = 74287464
$ws.Range("J126").Value = 100001870
$ws.Range("L126").Value = 300005610
$ws.Range("N126").Value = -300010550
$ws.Range("H136").Value = 2324.9697
$ws.Range("I136").Value = 2000.5
$ws.Range("J136").Value = 3530.1428
$ws.Range("K136").Value = 6001.5
$ws.Range("L136").Value = 10590.4284
$ws.Range("M136").Value = -3451.5
$ws.Range("N136").Value = -15690.4284

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4280324.5
$ws.Range("I4").Value = 2500368.8
$ws.Range("J4").Value = 17333332
$ws.Range("K4").Value = 7501106.399999999
$ws.Range("L4").Value = 51999996
$ws.Range("M4").Value = -7500994.399999999
$ws.Range("N4").Value = -52000220

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 3214.75
$ws.Range("I31").Value = 3214.75
$ws.Range("K31").Value = 3214.75
$ws.Range("M31").Value = -2922.75
$ws.Range("H37").Value = 3214.75
$ws.Range("I37").Value = 3214.75
$ws.Range("K37").Value = 3214.75
$ws.Range("M37").Value = -2937.75
$ws.Range("H41").Value = 9745.666999999999
$ws.Range("I41").Value = 9745.666999999999
$ws.Range("K41").Value = 9745.666999999999
$ws.Range("M41").Value = -9390.666999999999
$ws.Range("H122").Value = 1900.0892
$ws.Range("I122").Value = 1661.762
$ws.Range("K122").Value = 4985.286
$ws.Range("M122").Value = -2535.286
$ws.Range("H124").Value = 133326.33
$ws.Range("J124").Value = 133326.33
$ws.Range("L124").Value = 133326.33
$ws.Range("N124").Value = -143146.33
$ws.Range("H132").Value = 575101.8
$ws.Range("I132").Value = 7252.7144
$ws.Range("J132").Value = 1071969.8
$ws.Range("K132").Value = 21758.1432
$ws.Range("L132").Value = 3215909.4
$ws.Range("M132").Value = -19228.1432
$ws.Range("N132").Value = -3220969.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1274.1818
$ws.Range("I16").Value = 1118.8966
$ws.Range("K16").Value = 1118.8966
$ws.Range("M16").Value = -948.8966
$ws.Range("H43").Value = 4137000
$ws.Range("I43").Value = 290000
$ws.Range("K43").Value = 290000
$ws.Range("M43").Value = -289807
$ws.Range("H56").Value = 7999.3335
$ws.Range("J56").Value = 4999
$ws.Range("L56").Value = 4999
$ws.Range("H100").Value = 2953.4
$ws.Range("I100").Value = 2579.3333
$ws.Range("K100").Value = 2579.3333
$ws.Range("M100").Value = -2038.3333
$ws.Range("N56").Value = -6381

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("H32").Value = 10037.333
$ws.Range("I32").Value = 11056
$ws.Range("K32").Value = 11056
$ws.Range("M32").Value = -10739
$ws.Range("H51").Value = 2915.8333
$ws.Range("I51").Value = 2915.8333
$ws.Range("K51").Value = 2915.8333
$ws.Range("M51").Value = -2405.8333
$ws.Range("H61").Value = 16234.429
$ws.Range("I61").Value = 16234.429
$ws.Range("K61").Value = 16234.429
$ws.Range("M61").Value = -15942.429
$ws.Range("H122").Value = 2695.4614
$ws.Range("I122").Value = 2004.1
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6012.299999999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3562.299999999999
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 2629.1162
$ws.Range("I132").Value = 2293.739
$ws.Range("J132").Value = 3014.8
$ws.Range("K132").Value = 6881.217000000001
$ws.Range("L132").Value = 9044.400000000001
$ws.Range("M132").Value = -4351.217000000001
$ws.Range("N132").Value = -14104.4
$ws.Range("M2").ClearContents()
